$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts the existing row 23 (and
# everything below it, through the old row 94) down by one, producing the
# old row94 -> new row95 tail and growing the used range to A1:R95.
$ws.Rows.Item(23).EntireRow.Insert()

# Populate the newly inserted row 23 with the new data record.
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44838
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112012
$ws.Range("G23").Value = "Espinaca"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 6500
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 6700
$ws.Range("N23").Value = "$/cuna 10 kilos"
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 670
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = "Hortaliza"
